$d = $word.ActiveDocument

# Collapse a range to the very end of the document body content (just
# before the sectPr, which Word's Content range does not include) so the
# new paragraphs land after the existing "Negative: ..." paragraph.
$end = $d.Content
$end.Collapse(0)

# Build a WordprocessingML fragment (wrapped in the Flat-OPC package format
# Word's Range.InsertXML expects) containing:
#   1. a brand-new empty paragraph
#   2. the "People: ..." prompt paragraph, with an explicit <w:tab/> run
#      between the positive and negative prompt text (matching the layout
#      used by the existing "Negative: ..." paragraph already in the doc)
$peopleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r><w:t xml:space="preserve">People: [_________] </w:t></w:r>
            <w:r><w:t>high resolution, lots of detail, concept art, studio lighting, rule of thirds, anatomically correct, good hands, good skin, proportional, visible hands, visible face</w:t></w:r>
            <w:r><w:tab/><w:t xml:space="preserve">Negative: </w:t></w:r>
            <w:r><w:t>ugly, bad composition, weird hands, weird face, bad composition, muddy colors, bad pose, bad silhouette, hard to read, ugly face, smudgy face, inhuman</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$end.InsertXML($peopleXml)

$d.Save()
